$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update EffectiveDate (F) and PreviousExpDate (I) from 08152023 -> 09212023
$ws.Range("F2").Value = "09212023"
$ws.Range("I2").Value = "09212023"
$ws.Range("F3").Value = "09212023"
$ws.Range("I3").Value = "09212023"

# Update ConstYear (N) from numeric 2022 -> text "2023"
$ws.Range("N2").Value = "2023"
$ws.Range("N3").Value = "2023"

# Update selection to C12
$ws.Range("C12").Select()

$wb.Save()
